$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "zzz" pending-marker moves down one row (de Assis 2011 results came in):
# Row 6: J6 goes from pending ("zzz") to done (1); K6 becomes done (1) too.
# Row 7: I7 goes from done (1) to pending ("zzz"); J7 becomes pending ("zzz") too.
# Row 8: I8 goes from pending ("zzz") to done (1).
# Row 9: I9 goes from pending ("zzz") to done (1).

$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 1

$ws.Range("I7").Value = "zzz"
$ws.Range("J7").Value = "zzz"

$ws.Range("I8").Value = 1

$ws.Range("I9").Value = 1

$ws.Range("K7").Select()
